$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G4").Value = "2016-09-06 04:58:22"

$zhcn.Range("H4").Value = "2016-09-06 04:58:17"
$zhcn.Range("K4").Value = "2016-09-06 04:58:43"

$dede.Range("H4").Value = "2016-09-06 04:58:22"
$dede.Range("K4").Value = "2016-09-06 04:58:51"
